# Apply run/paragraph formatting changes described by the diff:
#  - Body paragraph 1 ("Tohle je muj super model."): rPr color=0000ff, sz=64, szCs=64
#  - Body paragraph 2 ("Seznam trid"):                rPr color=000080, sz=56, szCs=56
#  - Footer paragraph ("Text v zapati"):               pPr jc=right,  rPr color=000080
#  - Header paragraph ("Text v zahlavi"):               pPr jc=center, rPr i/iCs
#
# InsertXML is used for the run-level rPr so that "complex script" twin
# properties (szCs / iCs) - which have no dedicated COM setter exposed by
# this host - are still produced. InsertXML replaces the exact Range it is
# called on, so the Range is first shrunk to exclude the trailing paragraph
# mark, leaving <w:pPr>/w:pStyle completely untouched.

function Build-RPrXml($rPrInner, $text) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr>$rPrInner</w:rPr><w:t xml:space="preserve">$text</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

$d = $word.ActiveDocument

# --- Body paragraph 1: "Tohle je muj super model." ---
# (rebuild the target via Document.Range(start,end): a Range re-fetched
#  straight off Paragraphs(n) can, for the *last* paragraph of a story,
#  make InsertXML insert instead of replace; a plain Range(start,end)
#  never has that problem)
$r1 = $d.Paragraphs(1).Range
$t1 = $d.Range($r1.Start, $r1.End - 1)
$t1.InsertXML((Build-RPrXml '<w:color w:val="0000ff"/><w:sz w:val="64"/><w:szCs w:val="64"/>' "Tohle je můj super model."))

# --- Body paragraph 2: "Seznam trid" ---
$r2 = $d.Paragraphs(2).Range
$t2 = $d.Range($r2.Start, $r2.End - 1)
$t2.InsertXML((Build-RPrXml '<w:color w:val="000080"/><w:sz w:val="56"/><w:szCs w:val="56"/>' "Seznam tříd"))

# --- Footer paragraph: "Text v zapati" ---
$footerDoc = $d.Sections(1).Footers(1)
$fr = $footerDoc.Range
$fr.End = $fr.End - 1
$fr.InsertXML((Build-RPrXml '<w:color w:val="000080"/>' "Text v zápatí"))
$footerDoc.Range.Paragraphs(1).Alignment = 2   # wdAlignParagraphRight

# --- Header paragraph: "Text v zahlavi" ---
$headerDoc = $d.Sections(1).Headers(1)
$hr = $headerDoc.Range
$hr.End = $hr.End - 1
$hr.InsertXML((Build-RPrXml '<w:i/><w:iCs/>' "Text v záhlaví"))
$headerDoc.Range.Paragraphs(1).Alignment = 1   # wdAlignParagraphCenter
